# FIX update VM table
# Adds a new "UpdateVMTable" row to the Library_Formula sheet, mirroring the
# pattern of the existing rows (CREATE/MODIFY | LIB_EWS | <method name> | | String | String,String)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

$ws.Cells.Item(20, 1).Value = "CREATE/MODIFY"
$ws.Cells.Item(20, 2).Value = "LIB_EWS"
$ws.Cells.Item(20, 3).Value = "UpdateVMTable"
$ws.Cells.Item(20, 5).Value = "String"
$ws.Cells.Item(20, 6).Value = "String,String"

# Move the active selection like the authored edit did (cursor ends up on D21)
[void]$ws.Range("D21").Select()
